# Ordenanza 1594 - reformat header block, underline article labels,
# add a footer with a restarted page numbering sequence (2071-2130 /
# "Agregamos 2071 a 2130").

$d = $word.ActiveDocument

# --- Paragraph 1: "Yerba Buena, 08 de Noviembre de 2007" (right aligned date) ---
$p1 = $d.Paragraphs(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12   # 240 twips

# --- Paragraph 2: "ORDENANZA Nº 1594" (now bold, spaced before/after) ---
$p2 = $d.Paragraphs(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12  # 240 twips
$p2.Format.SpaceAfter = 18   # 360 twips
$p2.Range.Font.Bold = $true

# --- Paragraph 3: "EL CONCEJO DELIBERANTE..." (now bold, indented block) ---
$p3 = $d.Paragraphs(3)
$p3.Format.KeepWithNext = $true
$p3.Format.SpaceBefore = 18  # 360 twips
$p3.Format.SpaceAfter = 18   # 360 twips
$p3.Format.LeftIndent = 99.2   # 1984 twips
$p3.Format.RightIndent = 99.2  # 1984 twips
$p3.Range.Font.Bold = $true

# --- Paragraph 4: "ARTICULO PRIMERO: APRUEBASE ..." ---
$p4 = $d.Paragraphs(4)
$p4.Format.KeepWithNext = $true
$p4.Format.SpaceAfter = 6    # 120 twips
$p4.Format.Alignment = 0     # remove justify ("both")

# --- Paragraph 5: "ARTICULO SEGUNDO: COMUNIQUESE ..." ---
$p5 = $d.Paragraphs(5)
$p5.Format.KeepWithNext = $true
$p5.Format.SpaceAfter = 6    # 120 twips
$p5.Format.Alignment = 0     # remove justify ("both")

# Underline "ARTICULO PRIMERO" and the colon right after it, leaving the
# following space un-underlined (splits the old ": " run in two).
$rng = $d.Content
$rng.Find.Execute("ARTICULO PRIMERO", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$rng.Font.Underline = 1
$colon = $d.Range($rng.End, $rng.End + 1)
$colon.Font.Underline = 1

# Same treatment for "ARTICULO SEGUNDO".
$rng2 = $d.Content
$rng2.Find.Execute("ARTICULO SEGUNDO", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$rng2.Font.Underline = 1
$colon2 = $d.Range($rng2.End, $rng2.End + 1)
$colon2.Font.Underline = 1

# --- Section: footer with (initially empty) page-number field, and the
#     page numbering restarted at 1997. ---
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$ftr.PageNumbers.Add()
$ftr.PageNumbers.StartingNumber = 1997

$fp = $ftr.Range.Paragraphs(1)
$fp.Range.Delete()
$fp.Range.Style = "Piedepgina"
$fp.Range.Font.Name = "Book Antiqua"
$fp.Range.Font.Size = 10
$fp.Range.Font.Color = 8421504

$fstyle = $d.Styles("Piedepgina")
$fstyle.NameLocal = "footer"
$fstyle.UnhideWhenUsed = $true
